$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 27.07648166666667
$ws.Range("H2").Value = 81.229445
$ws.Range("I2").Value = 0.1943552322922666
$ws.Range("J2").Value = 0.1943552322922666
$ws.Range("M2").Value = 2.824532
$ws.Range("N2").Value = 8.473596000000001
$ws.Range("O2").Value = 0.01352821635215845
$ws.Range("P2").Value = 0.01352821635215845
$ws.Range("Q2").Value = 76.47838891491334
$ws.Range("R2").Value = 688.3055002342201
$ws.Range("S2").Value = 0.002629279631623795
$ws.Range("T2").Value = 0.002629279631623796

$ws.Range("G3").Value = 27.07648166666667
$ws.Range("H3").Value = 81.229445
$ws.Range("I3").Value = 0.1943552322922666
$ws.Range("J3").Value = 0.1943552322922666
$ws.Range("O3").Value = 0.5142441883127264
$ws.Range("P3").Value = 0.5142441883127264
$ws.Range("Q3").Value = 2907.150950815455
$ws.Range("R3").Value = 26164.35855733909
$ws.Range("S3").Value = 0.09994604867446805
$ws.Range("T3").Value = 0.09994604867446805

$ws.Range("G4").Value = 27.07648166666667
$ws.Range("H4").Value = 81.229445
$ws.Range("I4").Value = 0.1943552322922666
$ws.Range("J4").Value = 0.1943552322922666
$ws.Range("M4").Value = 28.70592233333333
$ws.Range("N4").Value = 86.117767
$ws.Range("O4").Value = 0.1374882380208794
$ws.Range("P4").Value = 0.1374882380208794
$ws.Range("Q4").Value = 777.2553797832572
$ws.Range("R4").Value = 6995.298418049315
$ws.Range("S4").Value = 0.02672155843800245
$ws.Range("T4").Value = 0.02672155843800246

$ws.Range("G5").Value = 27.07648166666667
$ws.Range("H5").Value = 81.229445
$ws.Range("I5").Value = 0.1943552322922666
$ws.Range("J5").Value = 0.1943552322922666
$ws.Range("M5").Value = 69.88962933333333
$ws.Range("N5").Value = 209.668888
$ws.Range("O5").Value = 0.3347393573142358
$ws.Range("P5").Value = 0.3347393573142358
$ws.Range("Q5").Value = 1892.365267334129
$ws.Range("R5").Value = 17031.28740600716
$ws.Range("S5").Value = 0.06505834554817233
$ws.Range("T5").Value = 0.06505834554817233

$ws.Range("I6").Value = 0.599012687336886
$ws.Range("J6").Value = 0.599012687336886
$ws.Range("M6").Value = 2.824532
$ws.Range("N6").Value = 8.473596000000001
$ws.Range("O6").Value = 0.01352821635215845
$ws.Range("P6").Value = 0.01352821635215845
$ws.Range("Q6").Value = 235.7102750813907
$ws.Range("R6").Value = 2121.392475732516
$ws.Range("S6").Value = 0.008103573231981237
$ws.Range("T6").Value = 0.008103573231981239

$ws.Range("I7").Value = 0.599012687336886
$ws.Range("J7").Value = 0.599012687336886
$ws.Range("O7").Value = 0.5142441883127264
$ws.Range("P7").Value = 0.5142441883127264
$ws.Range("S7").Value = 0.3080387931885819
$ws.Range("T7").Value = 0.3080387931885819

$ws.Range("I8").Value = 0.599012687336886
$ws.Range("J8").Value = 0.599012687336886
$ws.Range("M8").Value = 28.70592233333333
$ws.Range("N8").Value = 86.117767
$ws.Range("O8").Value = 0.1374882380208794
$ws.Range("P8").Value = 0.1374882380208794
$ws.Range("Q8").Value = 2395.540517740651
$ws.Range("R8").Value = 21559.86465966586
$ws.Range("S8").Value = 0.08235719893410037
$ws.Range("T8").Value = 0.08235719893410039

$ws.Range("I9").Value = 0.599012687336886
$ws.Range("J9").Value = 0.599012687336886
$ws.Range("M9").Value = 69.88962933333333
$ws.Range("N9").Value = 209.668888
$ws.Range("O9").Value = 0.3347393573142358
$ws.Range("P9").Value = 0.3347393573142358
$ws.Range("Q9").Value = 5832.365770859184
$ws.Range("R9").Value = 52491.29193773265
$ws.Range("S9").Value = 0.2005131219822225
$ws.Range("T9").Value = 0.2005131219822225

$ws.Range("G10").Value = 14.445417
$ws.Range("H10").Value = 43.336251
$ws.Range("I10").Value = 0.1036893349422856
$ws.Range("J10").Value = 0.1036893349422857
$ws.Range("M10").Value = 2.824532
$ws.Range("N10").Value = 8.473596000000001
$ws.Range("O10").Value = 0.01352821635215845
$ws.Range("P10").Value = 0.01352821635215845
$ws.Range("Q10").Value = 40.801542569844
$ws.Range("R10").Value = 367.213883128596
$ws.Range("S10").Value = 0.001402731756510663
$ws.Range("T10").Value = 0.001402731756510664

$ws.Range("G11").Value = 14.445417
$ws.Range("H11").Value = 43.336251
$ws.Range("I11").Value = 0.1036893349422856
$ws.Range("J11").Value = 0.1036893349422857
$ws.Range("O11").Value = 0.5142441883127264
$ws.Range("P11").Value = 0.5142441883127264
$ws.Range("Q11").Value = 1550.977275536318
$ws.Range("R11").Value = 13958.79547982686
$ws.Range("S11").Value = 0.0533216378840821
$ws.Range("T11").Value = 0.05332163788408211

$ws.Range("G12").Value = 14.445417
$ws.Range("H12").Value = 43.336251
$ws.Range("I12").Value = 0.1036893349422856
$ws.Range("J12").Value = 0.1036893349422857
$ws.Range("M12").Value = 28.70592233333333
$ws.Range("N12").Value = 86.117767
$ws.Range("O12").Value = 0.1374882380208794
$ws.Range("P12").Value = 0.1374882380208794
$ws.Range("Q12").Value = 414.669018474613
$ws.Range("R12").Value = 3732.021166271517
$ws.Range("S12").Value = 0.01425606396277165
$ws.Range("T12").Value = 0.01425606396277166

$ws.Range("G13").Value = 14.445417
$ws.Range("H13").Value = 43.336251
$ws.Range("I13").Value = 0.1036893349422856
$ws.Range("J13").Value = 0.1036893349422857
$ws.Range("M13").Value = 69.88962933333333
$ws.Range("N13").Value = 209.668888
$ws.Range("O13").Value = 0.3347393573142358
$ws.Range("P13").Value = 0.3347393573142358
$ws.Range("Q13").Value = 1009.584839695432
$ws.Range("R13").Value = 9086.263557258886
$ws.Range("S13").Value = 0.03470890133892122
$ws.Range("T13").Value = 0.03470890133892123

$ws.Range("G14").Value = 14.34140633333333
$ws.Range("H14").Value = 43.024219
$ws.Range("I14").Value = 0.1029427454285617
$ws.Range("J14").Value = 0.1029427454285617
$ws.Range("M14").Value = 2.824532
$ws.Range("N14").Value = 8.473596000000001
$ws.Range("O14").Value = 0.01352821635215845
$ws.Range("P14").Value = 0.01352821635215845
$ws.Range("Q14").Value = 40.50776111350267
$ws.Range("R14").Value = 364.569850021524
$ws.Range("S14").Value = 0.001392631732042752
$ws.Range("T14").Value = 0.001392631732042753

$ws.Range("G15").Value = 14.34140633333333
$ws.Range("H15").Value = 43.024219
$ws.Range("I15").Value = 0.1029427454285617
$ws.Range("J15").Value = 0.1029427454285617
$ws.Range("O15").Value = 0.5142441883127264
$ws.Range("P15").Value = 0.5142441883127264
$ws.Range("Q15").Value = 1539.809845727031
$ws.Range("R15").Value = 13858.28861154328
$ws.Range("S15").Value = 0.05293770856559432
$ws.Range("T15").Value = 0.05293770856559433

$ws.Range("G16").Value = 14.34140633333333
$ws.Range("H16").Value = 43.024219
$ws.Range("I16").Value = 0.1029427454285617
$ws.Range("J16").Value = 0.1029427454285617
$ws.Range("M16").Value = 28.70592233333333
$ws.Range("N16").Value = 86.117767
$ws.Range("O16").Value = 0.1374882380208794
$ws.Range("P16").Value = 0.1374882380208794
$ws.Range("Q16").Value = 411.6832963554414
$ws.Range("R16").Value = 3705.149667198973
$ws.Range("S16").Value = 0.01415341668600488
$ws.Range("T16").Value = 0.01415341668600488

$ws.Range("G17").Value = 14.34140633333333
$ws.Range("H17").Value = 43.024219
$ws.Range("I17").Value = 0.1029427454285617
$ws.Range("J17").Value = 0.1029427454285617
$ws.Range("M17").Value = 69.88962933333333
$ws.Range("N17").Value = 209.668888
$ws.Range("O17").Value = 0.3347393573142358
$ws.Range("P17").Value = 0.3347393573142358
$ws.Range("Q17").Value = 1002.315572755386
$ws.Range("R17").Value = 9020.840154798472
$ws.Range("S17").Value = 0.03445898844491972
$ws.Range("T17").Value = 0.03445898844491972

